$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule R30's "Integer min" (C1) value changes from 18 to 1.
$ws.Range("C10").Value = 1
